$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 340
$ws.Range("I55").Value = 274
$ws.Range("J55").Value = 395
$ws.Range("K55").Value = 274
$ws.Range("L55").Value = 395
$ws.Range("M55").Value = -60
$ws.Range("N55").Value = -823

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2133551.5
$ws.Range("I76").Value = 3349003
$ws.Range("J76").Value = 6511.25
$ws.Range("K76").Value = 3349003
$ws.Range("L76").Value = 6511.25
$ws.Range("M76").Value = -3348688
$ws.Range("N76").Value = -7141.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2133551.5
$ws.Range("I79").Value = 3349003
$ws.Range("J79").Value = 6511.25
$ws.Range("K79").Value = 3349003
$ws.Range("L79").Value = 6511.25
$ws.Range("M79").Value = -3347911
$ws.Range("N79").Value = -8695.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1360.7142
$ws.Range("I98").Value = 1378.75
$ws.Range("K98").Value = 1378.75
$ws.Range("M98").Value = 119.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 33335466
$ws.Range("I111").Value = 100000000
$ws.Range("K111").Value = 300000000
$ws.Range("M111").Value = -299996933

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11962.714
$ws.Range("I116").Value = 18724.166
$ws.Range("J116").Value = 6891.625
$ws.Range("K116").Value = 18724.166
$ws.Range("L116").Value = 6891.625
$ws.Range("M116").Value = -15282.166
$ws.Range("N116").Value = -13775.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1360.7142
$ws.Range("I122").Value = 1378.75
$ws.Range("K122").Value = 4136.25
$ws.Range("M122").Value = -1686.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2825.28
$ws.Range("J138").Value = 2882.04
$ws.Range("L138").Value = 8646.119999999999
$ws.Range("N138").Value = -18926.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1453899.5
$ws.Range("I2").Value = 1938342.2
$ws.Range("J2").Value = 571.5
$ws.Range("K2").Value = 1938342.2
$ws.Range("L2").Value = 571.5
$ws.Range("M2").Value = -1938229.2
$ws.Range("N2").Value = -797.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3398.1355
$ws.Range("I32").Value = 2371.3396
$ws.Range("J32").Value = 12468.167
$ws.Range("K32").Value = 2371.3396
$ws.Range("L32").Value = 12468.167
$ws.Range("M32").Value = -2084.3396
$ws.Range("N32").Value = -13042.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1707.4667
$ws.Range("J45").Value = 1854.2222
$ws.Range("L45").Value = 1854.2222
$ws.Range("N45").Value = -2608.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3735.9092
$ws.Range("J88").Value = 4599.7144
$ws.Range("L88").Value = 4599.7144
$ws.Range("N88").Value = -5411.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3735.9092
$ws.Range("J91").Value = 4599.7144
$ws.Range("L91").Value = 4599.7144
$ws.Range("N91").Value = -7407.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1453899.5
$ws.Range("I116").Value = 1938342.2
$ws.Range("J116").Value = 571.5
$ws.Range("K116").Value = 1938342.2
$ws.Range("L116").Value = 571.5
$ws.Range("M116").Value = -1936048.2
$ws.Range("N116").Value = -5159.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6615.3335
$ws.Range("I122").Value = 6615.3335
$ws.Range("K122").Value = 19846.0005
$ws.Range("M122").Value = -17396.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1741.2
$ws.Range("I132").Value = 1127.8334
$ws.Range("K132").Value = 3383.5002
$ws.Range("M132").Value = -853.5001999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1453899.5
$ws.Range("I3").Value = 1938342.2
$ws.Range("J3").Value = 571.5
$ws.Range("K3").Value = 1938342.2
$ws.Range("L3").Value = 571.5
$ws.Range("M3").Value = -1938228.2
$ws.Range("N3").Value = -799.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6800.8667
$ws.Range("I80").Value = 44.75
$ws.Range("J80").Value = 9257.637000000001
$ws.Range("K80").Value = 44.75
$ws.Range("L80").Value = 9257.637000000001
$ws.Range("M80").Value = 953.25
$ws.Range("N80").Value = -11253.637

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 6800.8667
$ws.Range("I83").Value = 44.75
$ws.Range("J83").Value = 9257.637000000001
$ws.Range("K83").Value = 223.75
$ws.Range("L83").Value = 46288.185
$ws.Range("M83").Value = 4768.25
$ws.Range("N83").Value = -56272.185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 107400.9
$ws.Range("I86").Value = 2024.4615
$ws.Range("K86").Value = 2024.4615
$ws.Range("M86").Value = -901.4614999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 107400.9
$ws.Range("I89").Value = 2024.4615
$ws.Range("K89").Value = 10122.3075
$ws.Range("M89").Value = -4506.307499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 749.5454999999999
$ws.Range("I16").Value = 724.5
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 724.5
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -437.5
$ws.Range("N16").Value = -1574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1660
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 11000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 749.5454999999999
$ws.Range("I113").Value = 724.5
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 724.5
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1445.5
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2358.4
$ws.Range("I103").Value = 1600.1666
$ws.Range("J103").Value = 3495.75
$ws.Range("K103").Value = 4800.4998
$ws.Range("L103").Value = 10487.25
$ws.Range("M103").Value = -3921.4998
$ws.Range("N103").Value = -12245.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12395.305
$ws.Range("J131").Value = 14410.136
$ws.Range("L131").Value = 43230.408
$ws.Range("N131").Value = -53310.408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11645.363
$ws.Range("I70").Value = 15614.143
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 15614.143
$ws.Range("L70").Value = 4700
$ws.Range("M70").Value = -15344.143
$ws.Range("N70").Value = -5240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11645.363
$ws.Range("I73").Value = 15614.143
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 15614.143
$ws.Range("L73").Value = 4700
$ws.Range("M73").Value = -14678.143
$ws.Range("N73").Value = -6572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1142.5
$ws.Range("J97").Value = 2055.5
$ws.Range("L97").Value = 2055.5
$ws.Range("N97").Value = -3047.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1482.1666
$ws.Range("I122").Value = 1142.8182
$ws.Range("J122").Value = 2015.4286
$ws.Range("K122").Value = 3428.4546
$ws.Range("L122").Value = 6046.2858
$ws.Range("M122").Value = -978.4546
$ws.Range("N122").Value = -10946.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1168105.9
$ws.Range("I132").Value = 1426524.2
$ws.Range("K132").Value = 4279572.6
$ws.Range("M132").Value = -4277042.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4929.1177
$ws.Range("I40").Value = 1345.1818
$ws.Range("J40").Value = 11499.667
$ws.Range("K40").Value = 1345.1818
$ws.Range("L40").Value = 11499.667
$ws.Range("M40").Value = -1209.1818
$ws.Range("N40").Value = -11771.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3386.125
$ws.Range("I68").Value = 3181.5
$ws.Range("K68").Value = 3181.5
$ws.Range("M68").Value = -2432.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3386.125
$ws.Range("I71").Value = 3181.5
$ws.Range("K71").Value = 15907.5
$ws.Range("M71").Value = -12163.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7400.5454
$ws.Range("I122").Value = 5175.125
$ws.Range("K122").Value = 15525.375
$ws.Range("M122").Value = -13075.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3485.8333
$ws.Range("I136").Value = 1678.2858
$ws.Range("J136").Value = 4636.091
$ws.Range("K136").Value = 5034.857400000001
$ws.Range("L136").Value = 13908.273
$ws.Range("M136").Value = -2484.857400000001
$ws.Range("N136").Value = -19008.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 777.17645
$ws.Range("I107").Value = 536.9091
$ws.Range("K107").Value = 1610.7273
$ws.Range("M107").Value = 309.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 40420
$ws.Range("J121").Value = 40420
$ws.Range("L121").Value = 40420
$ws.Range("N121").Value = -43914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2654.1155
$ws.Range("I132").Value = 2223.2666
$ws.Range("K132").Value = 6669.7998
$ws.Range("M132").Value = -4139.7998
